$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$caseIdCohortQuery = 'MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
WHERE demo.breed IN [''American Staffordshire Terrier'']
MATCH (c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '''') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '''') AS `Study Code` ,
        coalesce(s.clinical_study_type, '''') AS  `Study Type`,
        coalesce(demo.breed, '''') AS Breed ,
        coalesce(diag.disease_term, '''') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '''') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '''') AS Age ,
        coalesce(demo.sex, '''') AS Sex ,
        coalesce(demo.neutered_indicator, '''') AS `Neutered Status`,
        coalesce(demo.weight, '''') AS `Weight (kg)`,
        coalesce(diag.best_response, '''') AS `Response to Treatment`,
        coalesce(co.cohort_description, '''') AS `Cohort`'

$sampleIdQuery = 'MATCH (s:study)<-[*]-(c:case)<--(demo:demographic), (samp:sample)-->(c)<--(diag:diagnosis) WHERE demo.breed IN [''American Staffordshire Terrier''] WITH DISTINCT samp AS samp, c, demo, diag
RETURN  coalesce(samp.sample_id, '''') AS `Sample ID`, 
        coalesce(c.case_id, '''') AS `Case ID`, 
        coalesce(demo.breed,'''') AS Breed , 
        coalesce(diag.disease_term,'''') AS Diagnosis , 
        coalesce(samp.sample_site, '''') AS `Sample Site`,
        coalesce(samp.summarized_sample_type, '''') AS `Sample Type`,
        coalesce(samp.specific_sample_pathology, '''') AS `Pathology/Morphology`,
        coalesce(samp.tumor_grade, '''') AS `Tumor Grade`,
        coalesce(samp.sample_chronology, '''') AS `Sample Chronology`,
        coalesce(samp.percentage_tumor, '''') AS `Percentage Tumor`,
        coalesce(samp.necropsy_sample, '''') AS `Necropsy Sample`,
        coalesce(samp.sample_preservation, '''') AS `Sample Preservation`'

$statCountQuery = 'MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE demo.breed IN[''American Staffordshire Terrier'']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study'

$ws.Range("B1").Value = "StatQuery"

# Write the Sample ID query text first so it lands earlier in the shared-string
# table (matches how the source workbook's string pool is ordered), then fill
# in the rest of row 3 and row 2.
$ws.Range("A3").Value = $sampleIdQuery
$ws.Range("A2").Value = $caseIdCohortQuery

$ws.Range("B2").Value = $statCountQuery
$ws.Range("C2").Value = "TC02_Canine_Filter_Breed-AmerStaffd_Neo4jData.xlsx"
$ws.Range("D2").Value = "TC02_Canine_Filter_Breed-AmerStaffd_WebData.xlsx"

$ws.Range("B3").Value = $statCountQuery
$ws.Range("C3").Value = "TC02_Canine_Filter_Breed-AmerStaffd_Neo4jData.xlsx"
$ws.Range("D3").Value = "TC02_Canine_Filter_Breed-AmerStaffd_WebData.xlsx"

$ws.Range("A2").WrapText = $true
$ws.Range("A3").WrapText = $true
$ws.Range("B3").WrapText = $true

$ws.Rows.Item(2).RowHeight = 261
$ws.Rows.Item(3).RowHeight = 217.5

$ws.Columns.Item(4).ColumnWidth = 45.25651041666667

$ws.Range("B2").Select() | Out-Null
